$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values of rows 3, 5, 6, 7 for the columns that rotate:
# D (Fecha), K (Variedad), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg)
$cols = @("D", "K", "L", "M", "N", "O", "P", "S")
$rows = @(3, 5, 6, 7)

$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2()
    }
}

# Rotation: new row3 <- old row5, new row5 <- old row6, new row6 <- old row7, new row7 <- old row3
$mapping = @{ 3 = 5; 5 = 6; 6 = 7; 7 = 3 }

foreach ($r in $rows) {
    $srcRow = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $orig[$srcRow][$c]
    }
}
